$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new question rows (row 37 intentionally left blank, matching source)
$ws.Range("C36").Value = "Można utworzyć instancję klasy Klasa nazwa = new Klasa(), lub private Klasa nazwa. Jaka różnica"
$ws.Range("C38").Value = "lec67test komentarz"
$ws.Range("C39").Value = "lec80 main, jak przetestowac metodę readIntegers. Wejscie rozmiar tablicy, uzytkownik podaje inty a na wyjsciu tablica int[]"

# Copy style from an existing question cell (e.g. C23) onto the new cells so
# they keep the same font/number-format/wrap formatting used for this column
$ws.Range("C23").Copy()
$ws.Range("C36").PasteSpecial(-4122)
$ws.Range("C38").PasteSpecial(-4122)
$ws.Range("C39").PasteSpecial(-4122)

$ws.Rows.Item(36).RowHeight = 28.5
$ws.Rows.Item(39).RowHeight = 28.5

# Update the view to match the scrolled-down state recorded in the workbook
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("C40").Select()
